# Moved delivery organisation path to be before practitioner key on service contact.
#
# On the "Service Contacts" sheet, the column "delivery_organisation_path"
# (previously the last data column, R) is moved to sit right after
# episode_key, i.e. it becomes the new column D (pushing practitioner_key
# and everything after it one column to the right, up to funding_source
# which ends up in column R). The tags column (S) is unaffected.

$wb = $excel.ActiveWorkbook

$wsSC = $wb.Worksheets.Item("Service Contacts")

# Cut the "delivery_organisation_path" column (R) and insert it before
# column D ("practitioner_key"), shifting D:R right to E:S... (S itself,
# "service_contact_tags", stays put because the cut already removed R from
# the column sequence before the insert happens).
$wsSC.Range("R1").EntireColumn.Cut()
$wsSC.Range("D1").EntireColumn.Insert()

# The insert operation stamps explicit (and irrelevant) width metadata onto
# the columns that simply shifted over (E:P); clear that incidental
# formatting so only the genuinely-customized columns keep explicit widths.
$wsSC.Range("E1:P1").EntireColumn.ClearFormats()

# Update the saved selection/view state on Service Contacts to the full
# column D (matches the workbook being saved right after performing the move).
$wsSC.Range("D1:D1048576").Select()

# The K5 and K10+ reference sheets had their selection state updated too
# (the user had highlighted the whole delivery_organisation_path column,
# F, while working on the move above) - K5's selection was saved, then
# K10+'s selection was saved last so K10+ remains the active tab, matching
# the workbook's unchanged activeTab.
$wsK5 = $wb.Worksheets.Item("K5")
$wsK5.Range("F1:F5").Select()

$wsK10p = $wb.Worksheets.Item("K10+")
$wsK10p.Range("F1:F5").Select()
